# Add an "included" column (J) to the subjects sheet, marking whether
# each subject's data was included in the final analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("J1").Value = "included"

# Per-row inclusion flags, rows 2-22 (one per subject)
$included = @(1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1)

for ($i = 0; $i -lt $included.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $included[$i]
}

# Match the author's final selection on the sheet
$ws.Range("J14").Select()
